$wb = $excel.ActiveWorkbook

# Insert a new worksheet named "VERSION" positioned between "Staff Data" and "Sheet2"
$sheet2 = $wb.Worksheets.Item("Sheet2")
$versionSheet = $wb.Worksheets.Add($sheet2)
$versionSheet.Name = "VERSION"

# Populate the version info starting at row 6
$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

# Make the VERSION tab the active / selected sheet
$versionSheet.Select()
$versionSheet.Range("B6").Select()
